$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.366.51'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '3.498.31'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.74'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.33'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.62'
$ws.Range("E9").Value = '  +6.04%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.393'
$ws.Range("E11").Value = '  +4.12%  '
$ws.Range("D12").Value = '4.093.16'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000181'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '3.499.05'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '64.376.59'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.63'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.57'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.59'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("E22").Value = '  +2.97%  '
$ws.Range("D23").Value = '3.637.74'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.33'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("E27").Value = '  +2.68%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.19'
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("E33").Value = '  +5.86%  '
$ws.Range("D34").Value = '3.526.61'
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.94'
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '165.48'
$ws.Range("E40").Value = '  +2.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0787'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.76'
$ws.Range("E45").Value = '  -2.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.19'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.926'
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").Value = '2.404.05'
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("E51").Value = '  +0.04%  '
